# Refactor experiment parameters and add dynamic calculation for crossover
# and mutation rates. Rebuilds the experiment grid as a full factorial
# design over population size (B) x generations (C), computing:
#   Crossover (D) = ROUND(0.8 + 0.1 * LN(Population), 3)
#   Mutation  (E) = ROUND(1 / SQRT(Generations), 3)
# and fills in the simulated fitness / duration / balance results for each
# of the 16 experiment rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Deney No, Populasyon, Nesiller, Caprazlama, Mutasyon,
#             Fitness Agirliklari, Fitness Degeri, Toplam Sure, Makine Yuk
#             Dengesizligi, Tip Degisim Sayisi
$data = @(
    @(1, 50, 10, 1.191, 0.316, "(-2, -3, -10)", 6430.903783855611, 2661.513897768687, 6402.32860092934, 280),
    @(2, 50, 50, 1.191, 0.141, "(-2, -3, -10)", 6487.026814483518, 3446.449539288811, 6310.973671538432, 280),
    @(3, 50, 100, 1.191, 0.1, "(-2, -3, -10)", 6396.231945803867, 3273.517367814602, 6355.747133909991, 280),
    @(4, 50, 200, 1.191, 0.071, "(-2, -3, -10)", 6461.857192319964, 2545.72317992938, 6318.900305878747, 280),
    @(5, 100, 10, 1.261, 0.316, "(-2, -3, -10)", 6421.928887737818, 2752.292209571341, 6358.124382064607, 280),
    @(6, 100, 50, 1.261, 0.141, "(-2, -3, -10)", 6405.836051533373, 3579.368005061616, 6305.003425672308, 280),
    @(7, 100, 100, 1.261, 0.1, "(-2, -3, -10)", 6438.481335196971, 3001.695682825894, 6270.45400907426, 280),
    @(8, 100, 200, 1.261, 0.071, "(-2, -3, -10)", 6427.693220095725, 2578.315197740972, 6300.708048426569, 280),
    @(9, 150, 10, 1.301, 0.316, "(-2, -3, -10)", 6460.461252407366, 2974.980946830065, 6410.652579342721, 280),
    @(10, 150, 50, 1.301, 0.141, "(-2, -3, -10)", 6357.430453075702, 2962.33533298378, 6321.80880149837, 280),
    @(11, 150, 100, 1.301, 0.1, "(-2, -3, -10)", 6407.39492048641, 2810.527906002668, 6297.908681315195, 280),
    @(12, 150, 200, 1.301, 0.071, "(-2, -3, -10)", 6376.158185460904, 2926.56359130191, 6244.618389071275, 280),
    @(13, 200, 10, 1.33, 0.316, "(-2, -3, -10)", 6425.265875087587, 3008.761377258846, 6327.74590398671, 280),
    @(14, 200, 50, 1.33, 0.141, "(-2, -3, -10)", 6370.464144974343, 3299.360630169088, 6329.273835921167, 280),
    @(15, 200, 100, 1.33, 0.1, "(-2, -3, -10)", 6392.628228112985, 2768.444900400293, 6303.592792956696, 280),
    @(16, 200, 200, 1.33, 0.071, "(-2, -3, -10)", 6453.684345521685, 2848.611804193922, 6319.97599708213, 280)
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
    $ws.Cells.Item($row, 8).Value = $rec[7]
    $ws.Cells.Item($row, 9).Value = $rec[8]
    $ws.Cells.Item($row, 10).Value = $rec[9]
    $row = $row + 1
}
